# feat: add 2022-Q4 data
#
# The "总计" (summary) sheet gains a new top row for 2022-Q4 (the other
# rows just slide down, keeping the same per-row values they already had
# one row up) and a brand-new last row that repeats the old 2021-Q4 total.
#
# A new "2022-Q4" worksheet is inserted right after "总计" (and therefore
# right before the existing "2022-Q3" tab) carrying the new quarterly fund
# holdings. It is built by copying the "2022-Q3" sheet (so it inherits the
# exact same layout/styles) and then overwriting just the data cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" summary sheet.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B2").Value = "2022-Q4"
# D2 keeps its existing value (0.1) -- nothing to change there.

$total.Range("B3").Value = "2022-Q3"
$total.Range("D3").Value = 0.1

$total.Range("B4").Value = "2022-Q2"
$total.Range("D4").Value = 0.16

$total.Range("B5").Value = "2022-Q1"
$total.Range("D5").Value = 0.15

# New row 6, reusing row 5's formatting for column A (bold/bordered style).
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q4"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.04

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet before "2022-Q3", cloning that
#    sheet's layout/styling and then replacing the figures with the new
#    quarter's numbers (fund codes/names are unchanged).
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

$q4.Range("D2:G3").NumberFormat = "@"

$q4.Range("D2").Value = "1.89"
$q4.Range("E2").Value = "90.53"
$q4.Range("F2").Value = "3.65"
$q4.Range("G2").Value = "0.0690"
$q4.Range("H2").Value = 6

$q4.Range("D3").Value = "0.97"
$q4.Range("E3").Value = "90.53"
$q4.Range("F3").Value = "3.65"
$q4.Range("G3").Value = "0.0354"
$q4.Range("H3").Value = 6
